# Remove rows 3 and 4 (country_codes_raw, iso_language_codes) leaving only
# the iso_language_codes_raw feed row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:4").Delete()

# L2's target org changed from "yalsworld_datawarehouse" to "sdmf"
$ws.Range("L2").Value = "sdmf"

# Update data-validation lists that gained new allowed options / a renamed org
$ws.Range("I2").Validation.Modify(3, 1, 1, '"DELTA_TABLE,JSON,CSV,PARQUET,XLSX,XML"')
$ws.Range("K2").Validation.Modify(3, 1, 1, '"FULL_LOAD,APPEND_LOAD,INCREMENTAL_CDC,SCD_TYPE_2,API_EXTRACTOR,STORAGE_FETCH"')
$ws.Range("L2").Validation.Modify(3, 1, 1, '"testing,sdmf"')

# Shrink the conditional-formatting ranges that used to span rows 2-4 down to
# row 2 only, now that rows 3-4 are gone.
$fcs = $ws.Cells.FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A2"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("A2"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("A2:R2"))
$fcs.Item(4).ModifyAppliesToRange($ws.Range("H2"))
$fcs.Item(5).ModifyAppliesToRange($ws.Range("N2"))

# Update the selection left over in the sheet view.
$ws.Range("F15").Select()
